$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9167996666666666
$ws.Range("H2").Value = 2.750399
$ws.Range("I2").Value = 0.2833456974325495
$ws.Range("J2").Value = 0.2833456974325495
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8534766666666666
$ws.Range("N2").Value = 2.56043
$ws.Range("O2").Value = 0.01737434360524059
$ws.Range("P2").Value = 0.01737434360524059
$ws.Range("Q2").Value = 0.7824671235077776
$ws.Range("R2").Value = 7.042204111569998
$ws.Range("S2").Value = 0.004922945506259652
$ws.Range("T2").Value = 0.004922945506259652

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9167996666666666
$ws.Range("H3").Value = 2.750399
$ws.Range("I3").Value = 0.2833456974325495
$ws.Range("J3").Value = 0.2833456974325495
$ws.Range("O3").Value = 0.7084105963118495
$ws.Range("P3").Value = 0.7084105963118494
$ws.Range("Q3").Value = 31.90382406109244
$ws.Range("R3").Value = 287.134416549832
$ws.Range("S3").Value = 0.2007250944805893
$ws.Range("T3").Value = 0.2007250944805893

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9167996666666666
$ws.Range("H4").Value = 2.750399
$ws.Range("I4").Value = 0.2833456974325495
$ws.Range("J4").Value = 0.2833456974325495
$ws.Range("O4").Value = 0.2742150600829099
$ws.Range("P4").Value = 0.2742150600829099
$ws.Range("Q4").Value = 12.34948923312811
$ws.Range("R4").Value = 111.145403098153
$ws.Range("S4").Value = 0.07769765744570056
$ws.Range("T4").Value = 0.07769765744570056

$ws.Range("I5").Value = 0.2271242616180895
$ws.Range("J5").Value = 0.2271242616180895
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8534766666666666
$ws.Range("N5").Value = 2.56043
$ws.Range("O5").Value = 0.01737434360524059
$ws.Range("P5").Value = 0.01737434360524059
$ws.Range("Q5").Value = 0.6272100451055554
$ws.Range("R5").Value = 5.644890405949999
$ws.Range("S5").Value = 0.003946134962439244
$ws.Range("T5").Value = 0.003946134962439244

$ws.Range("I6").Value = 0.2271242616180895
$ws.Range("J6").Value = 0.2271242616180895
$ws.Range("O6").Value = 0.7084105963118495
$ws.Range("P6").Value = 0.7084105963118494
$ws.Range("S6").Value = 0.1608972336097593
$ws.Range("T6").Value = 0.1608972336097593

$ws.Range("I7").Value = 0.2271242616180895
$ws.Range("J7").Value = 0.2271242616180895
$ws.Range("O7").Value = 0.2742150600829099
$ws.Range("P7").Value = 0.2742150600829099
$ws.Range("S7").Value = 0.06228089304589095
$ws.Range("T7").Value = 0.06228089304589095

$ws.Range("G8").Value = 1.583934333333334
$ws.Range("H8").Value = 4.751803000000001
$ws.Range("I8").Value = 0.4895300409493609
$ws.Range("J8").Value = 0.4895300409493609
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8534766666666666
$ws.Range("N8").Value = 2.56043
$ws.Range("O8").Value = 0.01737434360524059
$ws.Range("P8").Value = 0.01737434360524059
$ws.Range("Q8").Value = 1.351850995032222
$ws.Range("R8").Value = 12.16665895529
$ws.Range("S8").Value = 0.008505263136541694
$ws.Range("T8").Value = 0.008505263136541694

$ws.Range("G9").Value = 1.583934333333334
$ws.Range("H9").Value = 4.751803000000001
$ws.Range("I9").Value = 0.4895300409493609
$ws.Range("J9").Value = 0.4895300409493609
$ws.Range("O9").Value = 0.7084105963118495
$ws.Range("P9").Value = 0.7084105963118494
$ws.Range("Q9").Value = 55.11952516161157
$ws.Range("R9").Value = 496.0757264545041
$ws.Range("S9").Value = 0.3467882682215009
$ws.Range("T9").Value = 0.3467882682215008

$ws.Range("G10").Value = 1.583934333333334
$ws.Range("H10").Value = 4.751803000000001
$ws.Range("I10").Value = 0.4895300409493609
$ws.Range("J10").Value = 0.4895300409493609
$ws.Range("O10").Value = 0.2742150600829099
$ws.Range("P10").Value = 0.2742150600829099
$ws.Range("S10").Value = 0.1342365095913183
$ws.Range("T10").Value = 0.1342365095913183
